$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(4, 0, 4, 2),
    @(3, 1, 3, 2),
    @(5, 1, 6, 2),
    @(5, 2, 4, 0),
    @(5, 1, 6, 2),
    @(4, 0, 4, 2),
    @(4, 2, 5, 1),
    @(4, 3, 2, 0),
    @(5, 1, 5, 2),
    @(5, 0, 2, 2),
    @(4, 0, 3, 3),
    @(5, 0, 4, 2),
    @(5, 2, 5, 1),
    @(6, 1, 6, 2),
    @(3, 3, 3, 0),
    @(6, 0, 5, 2),
    @(3, 3, 3, 0),
    @(4, 2, 6, 0),
    @(6, 0, 5, 2),
    @(4, 2, 3, 1)
)

$startRow = 1922
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}

$ws.Range("A1942").Select()